$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Tere tulemast..." intro text to row 2 (general_intro_1) with a
# small wording fix ("Professor" -> "professor"), and make sure the other
# rows keep their original content (they merely shift shared-string indices
# as a side effect, which Excel manages automatically on save).
$introText = "Tere tulemast ja aitäh, et nõustusite osalema meie eksperimendis!\n\n`nSee uuring on osa rahvusvahelisest algatusest #EEGManyLabs, mille eesmärk on uurida ligikaudu 20 mõjuka kognitiiv-neuroteadusliku EEG-uuringu korratavust. Selles eksperimendis kordame ja laiendame professor Matthias M. Mülleri ja tema kolleegide 2003. aastal läbi viidud ruumilise tähelepanu uuringut."

$ws.Range("B2").Value = $introText

# Restore the selection / view state to match the saved workbook: active
# cell B3 and no pinned top-left scroll cell.
$ws.Range("B3").Select()
